$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1071.3846
$ws.Range("J17").Value = 1071.3846
$ws.Range("L17").Value = 3214.1538
$ws.Range("N17").Value = -3550.1538
$ws.Range("H74").Value = 4482.5557
$ws.Range("I74").Value = 3667.875
$ws.Range("K74").Value = 3667.875
$ws.Range("M74").Value = -2731.875
$ws.Range("H76").Value = 3898
$ws.Range("I76").Value = 3898
$ws.Range("K76").Value = 3898
$ws.Range("M76").Value = -3583
$ws.Range("H77").Value = 4482.5557
$ws.Range("I77").Value = 3667.875
$ws.Range("K77").Value = 18339.375
$ws.Range("M77").Value = -13659.375
$ws.Range("H79").Value = 3898
$ws.Range("I79").Value = 3898
$ws.Range("K79").Value = 3898
$ws.Range("M79").Value = -2806
$ws.Range("H96").Value = 1339.5625
$ws.Range("J96").Value = 978.44446
$ws.Range("L96").Value = 2935.33338
$ws.Range("N96").Value = -5681.33338
$ws.Range("H131").Value = 1649.25
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H137").Value = 2328
$ws.Range("I137").Value = 2328
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6984
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -4434
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 3347.125
$ws.Range("J138").Value = 3578.2
$ws.Range("L138").Value = 10734.6
$ws.Range("N138").Value = -21014.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8744.833000000001
$ws.Range("I31").Value = 6493.8
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 6493.8
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -6199.8
$ws.Range("N31").Value = -20588
$ws.Range("H45").Value = 3159.2
$ws.Range("I45").Value = 2800
$ws.Range("K45").Value = 2800
$ws.Range("M45").Value = -2423
$ws.Range("H74").Value = 1851.5
$ws.Range("I74").Value = 1570.2727
$ws.Range("J74").Value = 2882.6667
$ws.Range("K74").Value = 1570.2727
$ws.Range("L74").Value = 2882.6667
$ws.Range("M74").Value = -696.2727
$ws.Range("N74").Value = -4630.6667
$ws.Range("H77").Value = 1851.5
$ws.Range("I77").Value = 1570.2727
$ws.Range("J77").Value = 2882.6667
$ws.Range("K77").Value = 7851.363499999999
$ws.Range("L77").Value = 14413.3335
$ws.Range("M77").Value = -3483.363499999999
$ws.Range("N77").Value = -23149.3335
$ws.Range("H88").Value = 2778.2222
$ws.Range("I88").Value = 2999
$ws.Range("J88").Value = 2750.625
$ws.Range("K88").Value = 2999
$ws.Range("L88").Value = 2750.625
$ws.Range("M88").Value = -2593
$ws.Range("N88").Value = -3562.625
$ws.Range("H91").Value = 2778.2222
$ws.Range("I91").Value = 2999
$ws.Range("J91").Value = 2750.625
$ws.Range("K91").Value = 2999
$ws.Range("L91").Value = 2750.625
$ws.Range("M91").Value = -1595
$ws.Range("N91").Value = -5558.625
$ws.Range("H102").Value = 1550
$ws.Range("I102").Value = 1550
$ws.Range("K102").Value = 1550
$ws.Range("M102").Value = 72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5433.5356
$ws.Range("I22").Value = 1414.2307
$ws.Range("J22").Value = 8916.933999999999
$ws.Range("K22").Value = 1414.2307
$ws.Range("L22").Value = 8916.933999999999
$ws.Range("M22").Value = -1064.2307
$ws.Range("N22").Value = -9616.933999999999
$ws.Range("H50").Value = 11484
$ws.Range("I50").Value = 21052
$ws.Range("K50").Value = 21052
$ws.Range("M50").Value = -20427
$ws.Range("H58").Value = 3539.8
$ws.Range("I58").Value = 2928
$ws.Range("J58").Value = 3947.6667
$ws.Range("K58").Value = 2928
$ws.Range("L58").Value = 3947.6667
$ws.Range("M58").Value = -2725
$ws.Range("N58").Value = -4353.6667
$ws.Range("H60").Value = 950
$ws.Range("I60").Value = 950
$ws.Range("K60").Value = 950
$ws.Range("M60").Value = -439
$ws.Range("H99").Value = 3500
$ws.Range("I99").Value = 3500
$ws.Range("K99").Value = 3500
$ws.Range("M99").Value = -2002
$ws.Range("H105").Value = 1168.3334
$ws.Range("I105").Value = 649.75
$ws.Range("K105").Value = 649.75
$ws.Range("M105").Value = 1097.25
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 3500
$ws.Range("K126").Value = 10500
$ws.Range("M126").Value = -8030
$ws.Range("H132").Value = 1705
$ws.Range("I132").Value = 1705
$ws.Range("K132").Value = 5115
$ws.Range("M132").Value = -2585
$ws.Range("H136").Value = 3539.8
$ws.Range("I136").Value = 2928
$ws.Range("J136").Value = 3947.6667
$ws.Range("K136").Value = 8784
$ws.Range("L136").Value = 11843.0001
$ws.Range("M136").Value = -6234
$ws.Range("N136").Value = -16943.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 48.57143
$ws.Range("I40").Value = 31.666666
$ws.Range("K40").Value = 126.666664
$ws.Range("M40").Value = -57.666664
$ws.Range("H137").Value = 4803.125
$ws.Range("I137").Value = 3299
$ws.Range("J137").Value = 5705.6
$ws.Range("K137").Value = 9897
$ws.Range("L137").Value = 17116.8
$ws.Range("M137").Value = -4797
$ws.Range("N137").Value = -27316.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 95000
$ws.Range("J45").Value = 95000
$ws.Range("L45").Value = 95000
$ws.Range("N45").Value = -96118
$ws.Range("H64").Value = 47500
$ws.Range("I64").Value = 45000
$ws.Range("K64").Value = 45000
$ws.Range("M64").Value = -44752
$ws.Range("H67").Value = 47500
$ws.Range("I67").Value = 45000
$ws.Range("K67").Value = 45000
$ws.Range("M67").Value = -44142
$ws.Range("H97").Value = 746.8570999999999
$ws.Range("J97").Value = 857.6
$ws.Range("L97").Value = 857.6
$ws.Range("N97").Value = -1849.6
$ws.Range("H102").Value = 5429.8
$ws.Range("I102").Value = 5429.8
$ws.Range("K102").Value = 5429.8
$ws.Range("M102").Value = -3807.8
$ws.Range("H122").Value = 3045.1538
$ws.Range("I122").Value = 2187.6667
$ws.Range("K122").Value = 6563.000100000001
$ws.Range("M122").Value = -4113.000100000001
$ws.Range("H126").Value = 2732.923
$ws.Range("I126").Value = 2732.923
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8198.769
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5728.769
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 499
$ws.Range("I7").Value = 499
$ws.Range("K7").Value = 499
$ws.Range("M7").Value = -387
$ws.Range("H22").Value = 2440.818
$ws.Range("I22").Value = 2300
$ws.Range("J22").Value = 2521.2856
$ws.Range("K22").Value = 2300
$ws.Range("L22").Value = 2521.2856
$ws.Range("M22").Value = -2005
$ws.Range("N22").Value = -3111.2856
$ws.Range("H27").Value = 2440.818
$ws.Range("I27").Value = 2300
$ws.Range("J27").Value = 2521.2856
$ws.Range("K27").Value = 2300
$ws.Range("L27").Value = 2521.2856
$ws.Range("M27").Value = -2193
$ws.Range("N27").Value = -2735.2856
$ws.Range("H40").Value = 3142.1333
$ws.Range("I40").Value = 3361
$ws.Range("K40").Value = 3361
$ws.Range("M40").Value = -3225
$ws.Range("H46").Value = 3062.375
$ws.Range("H93").Value = 1317
$ws.Range("I93").Value = 1317
$ws.Range("K93").Value = 1317
$ws.Range("M93").Value = -69
$ws.Range("H100").Value = 2610.2856
$ws.Range("I100").Value = 1878.6666
$ws.Range("K100").Value = 1878.6666
$ws.Range("M100").Value = -1337.6666
$ws.Range("H122").Value = 6717.1665
$ws.Range("J122").Value = 2500
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 499
$ws.Range("I126").Value = 499
$ws.Range("K126").Value = 1497
$ws.Range("M126").Value = 973

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4791.4614
$ws.Range("I81").Value = 4133
$ws.Range("J81").Value = 5355.857
$ws.Range("K81").Value = 8266
$ws.Range("L81").Value = 10711.714
$ws.Range("M81").Value = -7205
$ws.Range("N81").Value = -12833.714
$ws.Range("H84").Value = 4791.4614
$ws.Range("I84").Value = 4133
$ws.Range("J84").Value = 5355.857
$ws.Range("K84").Value = 41330
$ws.Range("L84").Value = 53558.57
$ws.Range("M84").Value = -36026
$ws.Range("N84").Value = -64166.57
$ws.Range("H100").Value = 1633.3334
$ws.Range("I100").Value = 450
$ws.Range("K100").Value = 900
$ws.Range("M100").Value = -359
$ws.Range("H107").Value = 2234.4666
$ws.Range("I107").Value = 1590.7778
$ws.Range("J107").Value = 3200
$ws.Range("K107").Value = 4772.3334
$ws.Range("L107").Value = 9600
$ws.Range("M107").Value = -2852.3334
$ws.Range("N107").Value = -13440
$ws.Range("H113").Value = 873
$ws.Range("I113").Value = 873
$ws.Range("K113").Value = 2619
$ws.Range("M113").Value = -449
$ws.Range("H122").Value = 7737
$ws.Range("I122").Value = 7737
$ws.Range("K122").Value = 23211
$ws.Range("M122").Value = -20761
$ws.Range("H132").Value = 2455.077
$ws.Range("I132").Value = 1067.3334
$ws.Range("J132").Value = 5577.5
$ws.Range("K132").Value = 3202.0002
$ws.Range("L132").Value = 16732.5
$ws.Range("M132").Value = -672.0001999999999
$ws.Range("N132").Value = -21792.5
$ws.Range("H136").Value = 3577.1875
$ws.Range("I136").Value = 3374.6875
$ws.Range("K136").Value = 10124.0625
$ws.Range("M136").Value = -7574.0625
